$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most replacement values are non-numeric-looking text (contain '*', '(', ')', etc.)
# and Excel's COM layer keeps them as text automatically. A handful of the new
# values are plain numeric-looking strings (e.g. "-0.08", "0.10"); assigning
# those straight to .Value would get auto-coerced to real numbers (and can even
# lose trailing zeros / introduce float noise), so for those specific cells we
# briefly force Text number-format, type the value, then clear the format again
# so the cell's style reverts to the sheet's default (matching how the rest of
# the table's text-but-numeric-looking cells, e.g. "41422", are stored).

$values = @{
    "B2" = "0.42***"
    "C2" = "-0.08"
    "D2" = "-0.05"
    "E2" = "0.64**"
    "F2" = "0.78***"
    "G2" = "0.85***"
    "B3" = "(0.13)"
    "C3" = "(0.14)"
    "D3" = "(0.15)"
    "E3" = "(0.29)"
    "F3" = "(0.28)"
    "G3" = "(0.29)"
    "B4" = "-0.20***"
    "C4" = "0.21***"
    "D4" = "0.24***"
    "E4" = "0.23***"
    "F4" = "0.09**"
    "G4" = "0.11**"
    "B5" = "(0.03)"
    "D5" = "(0.04)"
    "E5" = "(0.04)"
    "F5" = "(0.04)"
    "G5" = "(0.04)"
    "B6" = "0.30***"
    "C6" = "0.13***"
    "D6" = "0.10***"
    "E6" = "0.09***"
    "F6" = "0.09***"
    "G6" = "0.11***"
    "F10" = "-0.33***"
    "D12" = "0.08***"
    "E12" = "0.08***"
    "F12" = "0.08***"
    "G12" = "0.08***"
    "G18" = "0.01***"
    "G20" = "0.01***"
    "G22" = "0.17***"
    "E24" = "-0.11***"
    "F24" = "-0.17***"
    "G24" = "-0.13***"
    "E25" = "(0.03)"
    "F25" = "(0.03)"
    "G25" = "(0.03)"
    "E26" = "-0.25"
    "F26" = "-0.31*"
    "G26" = "-0.17"
    "E27" = "(0.19)"
    "F27" = "(0.19)"
    "G27" = "(0.19)"
    "D29" = "0.03"
    "F29" = "0.10"
}

$numericLooking = @("C2", "D2", "E26", "G26", "D29", "F29")

foreach ($cell in $numericLooking) {
    $ws.Range($cell).NumberFormat = "@"
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

foreach ($cell in $numericLooking) {
    $ws.Range($cell).ClearFormats()
}
